# modify best bound and best objective in the tables
# The "Best Objective" / "Best Bound" header cells in columns C and D of the
# two result tables (row 3 and row 13) were swapped: what used to be in C
# (Best Objective) is now in D, and what used to be in D (Best Bound) is now
# in C. Swap both the text and the cell formatting (borders) that go with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-CellContentAndFormat {
    param(
        [string]$Addr1,
        [string]$Addr2,
        [string]$TempAddr
    )

    $r1 = $ws.Range($Addr1)
    $r2 = $ws.Range($Addr2)
    $tmp = $ws.Range($TempAddr)

    # Remember the current values of both cells.
    $v1 = $r1.Value()
    $v2 = $r2.Value()

    # Stash r1's formatting in a scratch cell, give r1 r2's formatting, then
    # give r2 the formatting that used to belong to r1.
    $r1.Copy() | Out-Null
    $tmp.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $r2.Copy() | Out-Null
    $r1.PasteSpecial(-4122) | Out-Null    # xlPasteFormats

    $tmp.Copy() | Out-Null
    $r2.PasteSpecial(-4122) | Out-Null    # xlPasteFormats

    # Swap the values to match the swapped formatting.
    $r1.Value = $v2
    $r2.Value = $v1

    # Clean up the scratch cell used for the format swap.
    $tmp.Clear() | Out-Null
    $excel.CutCopyMode = 0
}

Swap-CellContentAndFormat "C3" "D3" "Z1"
Swap-CellContentAndFormat "C13" "D13" "Z1"

# Update the saved cursor/selection position on the sheet.
$ws.Range("F9").Select() | Out-Null
